$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$modelText = "MultiOutputRegressor(estimator=GridSearchCV(cv=5,`n" + `
             "                                            estimator=Pipeline(steps=[('model',`n" + `
             "                                                                       RandomForestRegressor())]),`n" + `
             "                                            param_grid={'model__max_depth': [3,`n" + `
             "                                                                             5,`n" + `
             "                                                                             7],`n" + `
             "                                                        'model__n_estimators': [50,`n" + `
             "                                                                                100,`n" + `
             "                                                                                150]},`n" + `
             "                                            scoring='neg_mean_squared_error'))"

$ws.Range("A1").Copy($ws.Range("F1"))
$ws.Range("F1").Value = "Modelo"

$ws.Range("B2").Value = 0.4843335627135847
$ws.Range("C2").Value = 0.9857738545202548
$ws.Range("D2").Value = 0.5144539743107902
$ws.Range("F2").Value = $modelText

$ws.Range("B3").Value = 0.09447949688590218
$ws.Range("C3").Value = 0.9986942677713554
$ws.Range("D3").Value = 0.2484978502691945
$ws.Range("F3").Value = $modelText

$ws.Range("B4").Value = 0.03688000465082599
$ws.Range("C4").Value = 0.9996198856558589
$ws.Range("D4").Value = 0.1345741075077279
$ws.Range("F4").Value = $modelText

$ws.Range("B5").Value = 0.1074361815468557
$ws.Range("C5").Value = 0.9993600333205098
$ws.Range("D5").Value = 0.2226613007526392
$ws.Range("F5").Value = $modelText

$ws.Rows("2:5").AutoFit()
